$wb = $excel.ActiveWorkbook

# --- Users sheet: rename "Coartney Williams" -> "Coartney Trone" ---
$users = $wb.Worksheets.Item("Users")
$users.Range("A2").Value = "Coartney Trone"

# --- RateSheetManagement sheet: "Matthew White" moves from Intern/Financial
#     Analyst 1's Staff Member (B9) to a note in D9, replaced by "Jason Chan" ---
$rsm = $wb.Worksheets.Item("RateSheetManagement")
$rsm.Range("D9").Value = $rsm.Range("B9").Text
$rsm.Range("B9").Value = "Jason Chan"

# --- DealTeamMembers sheet: clear the stray explicit "default font" style
#     that had been applied to A2 ---
$dtm = $wb.Worksheets.Item("DealTeamMembers")
$dtm.Range("A2").ClearFormats()

# --- Re-create the selection / active-sheet trail left by the editing
#     session (Users A2 touched first, ends with RateSheetManagement B9
#     selected and that sheet active/visible on save) ---
$users.Activate()
$users.Range("A2").Select()

$dtm.Activate()
$dtm.Range("H7").Select()

$updateTimer = $wb.Worksheets.Item("Update_Timer")
$updateTimer.Activate()
$updateTimer.Range("G17").Select()

$rsm.Activate()
$rsm.Range("B9").Select()
